$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update version / date text in the intro cell (A1) ---
$ws.Range("A1").Value = "Onderstaande checklist kan gebruikt worden voor het uitvoeren van een assessment tegen de Kwaliteitsaanpak ICTU Software Realisatie versie 1.1.288, 12-07-2018."

# --- 2. Reorder / relabel the M07 "Continuous delivery pipeline" checklist sub-items (B29:B36) ---
# Original order:
#   1. bouw van de software
#   2. kwaliteitscontroles
#   3. regressietests
#   4. performance tests
#   5. beveiligingstests
#   6. unit tests
#   7. installatie van de software
#   8. oplevering van het totale product, ...
# New order:
#   1. bouw van de software
#   2. unit tests
#   3. regressietests
#   4. kwaliteitscontroles
#   5. performancetests (*)
#   6. beveiligingstests (*)
#   7. installatie van de software
#   8. oplevering van het totale product, ...
$ws.Range("B30").Value = "2. unit tests"
$ws.Range("B32").Value = "4. kwaliteitscontroles"
$ws.Range("B33").Value = "5. performancetests (*)"
$ws.Range("B34").Value = "6. beveiligingstests (*)"

# --- 3. Update the M07 comment text (cell B28) ---
$m07Comment = $ws.Range("B28").Comment
$m07Text = "Continuous delivery pipeline (M07)`n`nEr is een geautomatiseerde continuous delivery pipeline die aantoonbaar correct werkt en ten minste de volgende activiteiten uitvoert:`n`n- bouw van de software,`n- unit tests,`n- regressietests,`n- kwaliteitscontroles,`n- performancetests (*),`n- beveiligingstests (*),`n- installatie van de software,`n- oplevering van het totale product, dus inclusief alle deliverables, in de vorm zoals bruikbaar voor en afgesproken met de opdrachtgever.`n`n(*) Idealerwijs voert de geautomatiseerde continuous delivery pipeline ook performance tests en beveiligingstests uit. Vanwege de doorlooptijden van tests (met name van duurtesten) en licenties van testtools is dat niet altijd haalbaar. In dat geval vinden de performance tests en beveiligingstests periodiek en zo vaak mogelijk plaats, bij voorkeur dagelijks.`nDe projectenorganisatie voorziet in mensen en hulpmiddelen, zodat projecten deze pipeline kunnen toepassen. Projecten zijn verantwoordelijk voor de correcte werking van de pipelin.`n`nRationale`n`nSoftware incrementeel opleveren (zie [M05: Iteratief en incrementeel ontwikkelproces](#iteratief-en-incrementeel-ontwikkelproces-m05-)) vereist dat de software frequent gebouwd, getest en opgeleverd kan worden. Om dit efficiënt en foutvrij te doen, dient het proces van bouwen, testen en opleveren geautomatiseerd te zijn; een continuous delivery pipeline faciliteert dit.`n`nICTU`n`nICTU gebruikt Jenkins of Team Foundation Server (TFS) als tool voor de implementatie van de continuous delivery pipeline. De ICTU Release Manager ondersteunt de laatste stap (oplevering van het totale product).`n"
$m07Comment.Text($m07Text)

# --- 4. Update the M26 comment text (cell B65) ---
$m26Comment = $ws.Range("B65").Comment
$m26Text = "Periodieke beoordeling informatiebeveiliging (M26)`n`nProjecten laten periodiek een beveiligingstest uitvoeren. De code wordt zowel geautomatiseerd, als handmatig onderzocht op veelvoorkomende kwetsbaarheden door een beveiligingsexpert van buiten het project. Bevindingen uit de beveiligingstest worden vastgelegd als onderdeel van de werkvoorraad voor het ontwikkelproces (zie [M05: Iteratief en incrementeel ontwikkelproces](#iteratief-en-incrementeel-ontwikkelproces-m05-)).`n`nRationale`n`nDoor het inschakelen van actuele, specifieke expertise wordt de kans vergroot dat eventuele kwetsbaarheden in de gerealiseerde software tijdig herkend worden.`n`nICTU`n`nSoftware wordt minimaal bij iedere grote release of tenminste twee keer per jaar onderworpen aan een beveiligingstest door beveiligingsexperts die ICTU daarvoor inhuurt. Op basis van documentatie en architectuurstudie, crystalbox security audits (broncodescan) en penetratieaudits beoordelen deze experts of de software voldoet aan de projectspecifieke niet-functionele eisen met betrekking tot beveiliging, of bekende kwetsbaarheden (OWASP) vermeden zijn en in hoeverre voldoende invulling gegeven is aan de normen vanuit die vanuit BIR en SSD gelden.`n`nIndien door de opdrachtgever gewenst kunnen securitytesten door een onafhankelijke derde partij worden uitgevoerd in een daarvoor door de opdrachtgever beschikbaar gestelde omgeving. Dit kan zowel incidenteel als structureel worden ingericht. Afspraken hierover worden bij voorkeur al in de voorbereidingsfase gemaakt.`n`nDe beveiligingstesten vinden plaats in aanvulling op de door tools uitgevoerde continue beveiligingsanalyse van de gerealiseerde software, zie [M16: Verplichte tools](#verplichte-tools-m16-). Bevindingen uit zowel een beveiligingstest als de continue analyse worden in Jira als issue - gemarkeerd als beveiligingsbugreport - vastgelegd op de backlog van het project.`n"
$m26Comment.Text($m26Text)
